# Rename the "Scanner" sheet to "Session"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scanner")
$ws.Name = "Session"

# Delete row 3 (the second data row) so only the header + one data row remain
$ws.Rows.Item(3).Delete()
